# Natmi following Dr Hou advice
# Rewrites the Npy -> Npy1r LR-pair table: rows 2-5 get updated stats
# (new sending/target cluster pairing with "ECs"/"FAPs" as additional
# target clusters) and rows 6-11 are appended so every sending cluster
# (M1, M2) is crossed with every target cluster (ECs, FAPs, M2, Neutro, sCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("M1", "Npy", "Npy1r", "ECs", 3, 1, 2.659013, 7.977039, 0.4672138103908383, 0.4672138103908383, 1, 0.5, 0.141749, 0.283498, 0.008311773795175947, 0.007810267403127276, 0.376912433737, 2.261474602422, 0.003883375505950873, 0.003649064793586452),
  @("M1", "Npy", "Npy1r", "FAPs", 3, 1, 2.659013, 7.977039, 0.4672138103908383, 0.4672138103908383, 1, 0.3333333333333333, 0.035039, 0.105117, 0.002054591157674269, 0.002895935345626882, 0.093169156507, 0.838522408563, 0.0009599333635723189, 0.001353020987475845),
  @("M1", "Npy", "Npy1r", "M2", 3, 1, 2.659013, 7.977039, 0.4672138103908383, 0.4672138103908383, 1, 0.3333333333333333, 0.050715, 0.152145, 0.002973788936940283, 0.004191539742956915, 0.134851844295, 1.213666598655, 0.00138939526052599, 0.001958345254711535),
  @("M1", "Npy", "Npy1r", "Neutro", 3, 1, 2.659013, 7.977039, 0.4672138103908383, 0.4672138103908383, 3, 1, 2.104360666666667, 6.313082, 0.1233939558289581, 0.1739231266459359, 5.595522369355333, 50.35970132419799, 0.05765136028204631, 0.08125928671533605),
  @("M1", "Npy", "Npy1r", "sCs", 3, 1, 2.659013, 7.977039, 0.4672138103908383, 0.4672138103908383, 2, 1, 14.7221375, 29.444275, 0.8632658902812514, 0.811179130862353, 39.1463550002875, 234.878130001725, 0.4033297459787428, 0.3789940926397284),
  @("M2", "Npy", "Npy1r", "ECs", 3, 1, 3.032199333333333, 9.096598, 0.5327861896091618, 0.5327861896091618, 1, 0.5, 0.141749, 0.283498, 0.008311773795175947, 0.007810267403127276, 0.4298112233006667, 2.578867339804, 0.004428398289225075, 0.004161202609540825),
  @("M2", "Npy", "Npy1r", "FAPs", 3, 1, 3.032199333333333, 9.096598, 0.5327861896091618, 0.5327861896091618, 1, 0.3333333333333333, 0.035039, 0.105117, 0.002054591157674269, 0.002895935345626882, 0.1062452324406667, 0.9562070919660001, 0.00109465779410195, 0.001542914358151038),
  @("M2", "Npy", "Npy1r", "M2", 3, 1, 3.032199333333333, 9.096598, 0.5327861896091618, 0.5327861896091618, 1, 0.3333333333333333, 0.050715, 0.152145, 0.002973788936940283, 0.004191539742956915, 0.15377798919, 1.38400190271, 0.001584393676414293, 0.00223319448824538),
  @("M2", "Npy", "Npy1r", "Neutro", 3, 1, 3.032199333333333, 9.096598, 0.5327861896091618, 0.5327861896091618, 3, 1, 2.104360666666667, 6.313082, 0.1233939558289581, 0.1739231266459359, 6.380841010559555, 57.427569095036, 0.06574259554691181, 0.09266383993059987),
  @("M2", "Npy", "Npy1r", "sCs", 3, 1, 3.032199333333333, 9.096598, 0.5327861896091618, 0.5327861896091618, 2, 1, 14.7221375, 29.444275, 0.8632658902812514, 0.811179130862353, 44.64045551274167, 267.84273307645, 0.4599361443025087, 0.4321850382226247)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}